$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Traditional Trade")

# Insert a new column before column S (19th column), shifting all data
# from S onward one column to the right.
$ws.Range("S1").EntireColumn.Insert()

# Populate the header of the newly inserted column.
$ws.Range("S1").Value = "Sub brand"

# The autofilter range needs to grow by one column (A1:AN33 -> A1:AO33).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:AO33").AutoFilter()

# Defined names that pointed at the old filter range need to follow the
# widened range as well.
foreach ($n in $wb.Names) {
    $n.RefersTo = "='Traditional Trade'!`$A`$1:`$AO`$33"
}

# Restore/refresh the active-cell selection in the frozen (bottom-left) pane.
[void]$ws.Range("S2").Select()
